$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: headers (renamed from slug-style to human-readable "Municipio ..." labels,
# and reordered across columns)
$ws.Range("A1").Value = "Municipio superficie, medida"
$ws.Range("B1").Value = "Municipio zona desfavorecida"
$ws.Range("C1").Value = "Municipio montaña"
$ws.Range("D1").Value = "Municipio código"
$ws.Range("E1").Value = "Municipio nombre"

# Row 2: measure/dimension annotations
$ws.Range("A2").Value = "iaest-measure:municipio-superficie-medida"
$ws.Range("B2").Value = "iaest-measure:municipio-zona-desfavorecida"
$ws.Range("C2").Value = "iaest-measure:municipio-montana"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "sdmx-dimension:refArea"

# Row 3: medida/dim marker
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "dim"

# Row 4: xsd type / codelist reference
$ws.Range("A4").Value = "xsd:double"
$ws.Range("B4").Value = "xsd:string"
$ws.Range("C4").Value = "xsd:string"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "URI-Municipio"
